$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 83827
$ws.Cells.Item(2, 2).Value = "Ana Vitória Albuquerque"
$ws.Cells.Item(2, 4).Value = "Consulta medica"
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 45100
$ws.Cells.Item(2, 7).Value = 4043.7

# Row 3
$ws.Cells.Item(3, 1).Value = 16436
$ws.Cells.Item(3, 2).Value = "Mirella das Neves"
$ws.Cells.Item(3, 4).Value = "Outros"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 45100
$ws.Cells.Item(3, 7).Value = 2974.16

# Row 4
$ws.Cells.Item(4, 1).Value = 79120
$ws.Cells.Item(4, 2).Value = "Laura Novais"
$ws.Cells.Item(4, 3).Value = "Juridico"
$ws.Cells.Item(4, 4).Value = "Problemas pessoais"
$ws.Cells.Item(4, 6).Value = 45102
$ws.Cells.Item(4, 7).Value = 2682.13

# Row 5
$ws.Cells.Item(5, 1).Value = 89688
$ws.Cells.Item(5, 2).Value = "Davi Lucca Rezende"
$ws.Cells.Item(5, 4).Value = "Consulta medica"
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 45081
$ws.Cells.Item(5, 7).Value = 2555.15

# Row 6
$ws.Cells.Item(6, 1).Value = 10219
$ws.Cells.Item(6, 2).Value = "Sr. João Camargo"
$ws.Cells.Item(6, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(6, 4).Value = "Viagem de negocios"
$ws.Cells.Item(6, 6).Value = 45100
$ws.Cells.Item(6, 7).Value = 2889.28

# Row 7
$ws.Cells.Item(7, 1).Value = 50938
$ws.Cells.Item(7, 2).Value = "Ravy Rocha"
$ws.Cells.Item(7, 3).Value = "Marketing"
$ws.Cells.Item(7, 4).Value = "Viagem de negocios"
$ws.Cells.Item(7, 6).Value = 45095
$ws.Cells.Item(7, 7).Value = 6631.46

# Row 8
$ws.Cells.Item(8, 1).Value = 68766
$ws.Cells.Item(8, 2).Value = "Renan da Cunha"
$ws.Cells.Item(8, 3).Value = "Financeiro"
$ws.Cells.Item(8, 4).Value = "Viagem de negocios"
$ws.Cells.Item(8, 6).Value = 45106
$ws.Cells.Item(8, 7).Value = 4610.21

# Row 9
$ws.Cells.Item(9, 1).Value = 46567
$ws.Cells.Item(9, 2).Value = "Melissa Sampaio"
$ws.Cells.Item(9, 3).Value = "P&D"
$ws.Cells.Item(9, 4).Value = "Outros"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 45090
$ws.Cells.Item(9, 7).Value = 4489.55

# Row 10
$ws.Cells.Item(10, 1).Value = 3962
$ws.Cells.Item(10, 2).Value = "Fernando Fernandes"
$ws.Cells.Item(10, 3).Value = "Engenharia"
$ws.Cells.Item(10, 4).Value = "Problemas pessoais"
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 45078
$ws.Cells.Item(10, 7).Value = 7191.76

# Row 11
$ws.Cells.Item(11, 1).Value = 22791
$ws.Cells.Item(11, 2).Value = "Isaque Brito"
$ws.Cells.Item(11, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(11, 4).Value = "Outros"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 45086
$ws.Cells.Item(11, 7).Value = 7248.34
